$d = $word.ActiveDocument

# 1) Insert a new paragraph "New third sentence" right after "Second sentence"
#    and before "A list:".
$secondSentence = $d.Paragraphs(2)
$secondSentence.Range.InsertParagraphAfter()
$newThird = $d.Paragraphs(3)
$newThird.Range.Text = "New third sentence"

# 2) Remove the "4" run from the 4th bullet item, leaving the (now empty)
#    list paragraph itself in place.
$find = $d.Content.Find
$find.Text = "4"
$find.Forward = $true
$find.Wrap = 1
$find.Execute()
$four = $find.Parent
$four.Text = ""

# 3) Remove "Aquí iban las " (together with its proofErr spell-check wrapper
#    elements) from the paragraph that used to read "Aquí iban las
#    mayúsculas", leaving just "mayúsculas" (still spell-marked) behind.
$find2 = $d.Content.Find
$find2.Text = "mayúsculas"
$find2.Forward = $true
$find2.Wrap = 1
$find2.Execute()
$mayParaRange = $find2.Parent.Paragraphs(1).Range
$mayXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:r><w:t>mayúsculas</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
$mayParaRange.InsertXML($mayXml)

# 4) Drop the trailing period on "Y un poquito más."
$d.Content.Find.Execute("Y un poquito más.", $true, $false, $false, $false,
                         $false, $true, 1, $false, "Y un poquito más", 2)

# 5) Add a new paragraph "Pero sin el punto" right after it, inheriting the
#    es-ES language formatting of the paragraph it follows.
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$newLast = $d.Paragraphs($d.Paragraphs.Count)
$newLast.Range.Text = "Pero sin el punto"
